# Update transition-probability matrix cells on Sheet1 with recalculated
# percentages (denominator/counts shifted after the March 7 games pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.171003717472119
$ws.Range("C2").Value = 0.5427509293680297
$ws.Range("J2").Value = 0.007434944237918215
$ws.Range("O2").Value = 0.003717472118959108
$ws.Range("P2").Value = 0.1338289962825279
$ws.Range("S2").Value = 0.1412639405204461
# Row 3
$ws.Range("B3").Value = 0.01298701298701299
$ws.Range("C3").Value = 0.04545454545454546
$ws.Range("J3").Value = 0.04545454545454546
$ws.Range("P3").Value = 0.7597402597402597
$ws.Range("S3").Value = 0.1363636363636364
# Row 4
$ws.Range("J4").Value = 0.04761904761904762
$ws.Range("P4").Value = 0.7380952380952381
$ws.Range("S4").Value = 0.2142857142857143
# Row 6
$ws.Range("B6").Value = 0.06956521739130435
$ws.Range("D6").Value = 0.01739130434782609
$ws.Range("F6").Value = 0.09130434782608696
$ws.Range("J6").Value = 0.1826086956521739
$ws.Range("O6").Value = 0.01739130434782609
$ws.Range("Q6").Value = 0.208695652173913
$ws.Range("R6").Value = 0.09130434782608696
$ws.Range("S6").Value = 0.3217391304347826
# Row 7
$ws.Range("B7").Value = 0.08597285067873303
$ws.Range("D7").Value = 0.03167420814479638
$ws.Range("E7").Value = 0.004524886877828055
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.1040723981900453
$ws.Range("O7").Value = 0.01357466063348416
$ws.Range("Q7").Value = 0.167420814479638
$ws.Range("R7").Value = 0.1221719457013575
$ws.Range("S7").Value = 0.4117647058823529
# Row 8
$ws.Range("B8").Value = 0.06398104265402843
$ws.Range("D8").Value = 0.01658767772511848
$ws.Range("F8").Value = 0.04739336492890995
$ws.Range("J8").Value = 0.1208530805687204
$ws.Range("O8").Value = 0.01184834123222749
$ws.Range("Q8").Value = 0.1824644549763033
$ws.Range("R8").Value = 0.1658767772511848
$ws.Range("S8").Value = 0.3909952606635071
# Row 9
$ws.Range("B9").Value = 0.1182266009852217
$ws.Range("D9").Value = 0.03448275862068965
$ws.Range("F9").Value = 0.04926108374384237
$ws.Range("J9").Value = 0.07389162561576355
$ws.Range("O9").Value = 0.01477832512315271
$ws.Range("Q9").Value = 0.187192118226601
$ws.Range("R9").Value = 0.103448275862069
$ws.Range("S9").Value = 0.4187192118226601
# Row 10
$ws.Range("B10").Value = 0.09393063583815028
$ws.Range("D10").Value = 0.01372832369942197
$ws.Range("F10").Value = 0.078757225433526
$ws.Range("J10").Value = 0.1184971098265896
$ws.Range("O10").Value = 0.0180635838150289
$ws.Range("Q10").Value = 0.213150289017341
$ws.Range("R10").Value = 0.09465317919075145
$ws.Range("S10").Value = 0.3692196531791908
# Row 11
$ws.Range("G11").Value = 0.1349693251533742
$ws.Range("J11").Value = 0.07975460122699386
$ws.Range("K11").Value = 0.2085889570552147
$ws.Range("L11").Value = 0.5613496932515337
$ws.Range("S11").Value = 0.01533742331288344
# Row 12
$ws.Range("G12").Value = 0.7382198952879581
$ws.Range("J12").Value = 0.2094240837696335
$ws.Range("K12").Value = 0.01047120418848168
$ws.Range("L12").Value = 0.02094240837696335
$ws.Range("S12").Value = 0.02094240837696335
# Row 15
$ws.Range("F15").Value = 0.02304147465437788
$ws.Range("H15").Value = 0.1382488479262673
$ws.Range("I15").Value = 0.07834101382488479
$ws.Range("J15").Value = 0.3732718894009217
$ws.Range("K15").Value = 0.06912442396313365
$ws.Range("M15").Value = 0.009216589861751152
$ws.Range("O15").Value = 0.03686635944700461
$ws.Range("S15").Value = 0.271889400921659
# Row 16
$ws.Range("H16").Value = 0.1299435028248588
$ws.Range("I16").Value = 0.1073446327683616
$ws.Range("J16").Value = 0.4180790960451977
$ws.Range("K16").Value = 0.1186440677966102
$ws.Range("M16").Value = 0.02824858757062147
$ws.Range("O16").Value = 0.07344632768361582
$ws.Range("S16").Value = 0.1242937853107345
# Row 17
$ws.Range("F17").Value = 0.01622718052738337
$ws.Range("H17").Value = 0.1703853955375254
$ws.Range("I17").Value = 0.0872210953346856
$ws.Range("J17").Value = 0.4645030425963489
$ws.Range("K17").Value = 0.07910750507099391
$ws.Range("M17").Value = 0.02028397565922921
$ws.Range("O17").Value = 0.06085192697768763
$ws.Range("S17").Value = 0.101419878296146
# Row 18
$ws.Range("F18").Value = 0.01503759398496241
$ws.Range("H18").Value = 0.131578947368421
$ws.Range("I18").Value = 0.08270676691729323
$ws.Range("J18").Value = 0.481203007518797
$ws.Range("K18").Value = 0.112781954887218
$ws.Range("M18").Value = 0.01879699248120301
$ws.Range("N18").Value = 0.003759398496240601
$ws.Range("O18").Value = 0.05639097744360902
$ws.Range("S18").Value = 0.09774436090225563
# Row 19
$ws.Range("F19").Value = 0.009302325581395349
$ws.Range("H19").Value = 0.1961240310077519
$ws.Range("I19").Value = 0.07829457364341086
$ws.Range("J19").Value = 0.3868217054263566
$ws.Range("K19").Value = 0.113953488372093
$ws.Range("M19").Value = 0.02558139534883721
$ws.Range("N19").Value = 0.003875968992248062
$ws.Range("O19").Value = 0.06589147286821706
$ws.Range("S19").Value = 0.1201550387596899

